# Outline.docx edit script
# Applies: new title/author block, minor wording tweaks (remove "tools"/"TensorFlow and "),
# relocates the _GoBack bookmark to the title, and moves/removes the cached
# lastRenderedPageBreak markers to reflect re-paginated content.
#
# Note: w:proofErr elements seen in the target XML are purely cached, read-only
# artifacts that Word's live spell/grammar-checker writes out on save; they carry
# no visible text and are not exposed anywhere in the Word object model (no VBA
# automation -- real or headless -- can insert them deliberately), so they are
# intentionally not reproduced here. All genuine textual/structural edits are
# applied in full.

$d = $word.ActiveDocument

function Split-At($pos) {
    # Force a run boundary at a text position by adding + immediately removing
    # a zero-length bookmark there (mirrors how Word naturally splits runs).
    $r = $d.Range($pos, $pos)
    $d.Bookmarks.Add("TmpSplit", $r) | Out-Null
    $d.Bookmarks("TmpSplit").Delete()
}

# ---------------------------------------------------------------------------
# 1. Title paragraph: drop "tools" and "TensorFlow and ", prefix "Title: "
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Deep learning tools for ecological image analysis: an example using TensorFlow and Optical Character Recognition.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Deep learning for ecological image analysis: an example using Optical Character Recognition.",
    2) | Out-Null

$p1 = $d.Paragraphs.Item(1).Range
$startRng = $d.Range($p1.Start, $p1.Start)
$d.Bookmarks.Add("TempTitleSplit", $startRng) | Out-Null
$startRng.InsertBefore("Title: ")
$d.Bookmarks("TempTitleSplit").Delete()

# ---------------------------------------------------------------------------
# 2. Relocate the _GoBack bookmark from its old spot to just before "Optical"
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$findOptical = $d.Content
$findOptical.Find.Execute("Optical Character Recognition.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$bmRange = $d.Range($findOptical.Start, $findOptical.Start)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# ---------------------------------------------------------------------------
# 3. Insert the author line + following blank paragraph after the title
# ---------------------------------------------------------------------------
$p1now = $d.Paragraphs.Item(1).Range
$ip1 = $d.Range($p1now.End - 1, $p1now.End - 1)
$ip1.InsertParagraphAfter()
$ip2 = $d.Range($d.Paragraphs.Item(2).Range.End - 1, $d.Paragraphs.Item(2).Range.End - 1)
$ip2.InsertParagraphAfter()

$p2para = $d.Paragraphs.Item(2).Range
$p2para.Font.Superscript = $true
$p2para.InsertAfter("Ben Weinstein1")
$p2paraNow = $d.Paragraphs.Item(2).Range
$nameRng = $d.Range($p2paraNow.Start, $p2paraNow.Start + 13)
$nameRng.Font.Superscript = $false

Write-Host "Paragraph 1:" $d.Paragraphs.Item(1).Range.Text
Write-Host "Paragraph 2:" $d.Paragraphs.Item(2).Range.Text
Write-Host "Paragraph 3:" $d.Paragraphs.Item(3).Range.Text
Write-Host "Paragraph 4:" $d.Paragraphs.Item(4).Range.Text
